$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.2019543973941368
$ws.Range("C2").Value = 0.5570032573289903
$ws.Range("J2").Value = 0.01628664495114007
$ws.Range("P2").Value = 0.1302931596091205
$ws.Range("S2").Value = 0.09446254071661238
$ws.Range("B3").Value = 0.00546448087431694
$ws.Range("C3").Value = 0.0546448087431694
$ws.Range("J3").Value = 0.02185792349726776
$ws.Range("P3").Value = 0.726775956284153
$ws.Range("S3").Value = 0.1912568306010929
$ws.Range("J4").Value = 0.04
$ws.Range("P4").Value = 0.72
$ws.Range("S4").Value = 0.24
$ws.Range("B6").Value = 0.0481283422459893
$ws.Range("D6").Value = 0.0053475935828877
$ws.Range("F6").Value = 0.0427807486631016
$ws.Range("J6").Value = 0.2887700534759358
$ws.Range("O6").Value = 0.0160427807486631
$ws.Range("Q6").Value = 0.1443850267379679
$ws.Range("R6").Value = 0.1122994652406417
$ws.Range("S6").Value = 0.3422459893048128
$ws.Range("B7").Value = 0.07053941908713693
$ws.Range("D7").Value = 0.01659751037344398
$ws.Range("F7").Value = 0.05394190871369295
$ws.Range("J7").Value = 0.1203319502074689
$ws.Range("O7").Value = 0.01244813278008299
$ws.Range("Q7").Value = 0.2116182572614108
$ws.Range("R7").Value = 0.1203319502074689
$ws.Range("S7").Value = 0.3941908713692946
$ws.Range("B8").Value = 0.08378870673952642
$ws.Range("D8").Value = 0.02003642987249545
$ws.Range("F8").Value = 0.03460837887067395
$ws.Range("J8").Value = 0.1402550091074681
$ws.Range("O8").Value = 0.01457194899817851
$ws.Range("Q8").Value = 0.1985428051001822
$ws.Range("R8").Value = 0.08014571948998178
$ws.Range("S8").Value = 0.4280510018214936
$ws.Range("B9").Value = 0.1242937853107345
$ws.Range("D9").Value = 0.01694915254237288
$ws.Range("F9").Value = 0.02824858757062147
$ws.Range("J9").Value = 0.1129943502824859
$ws.Range("Q9").Value = 0.1864406779661017
$ws.Range("R9").Value = 0.1016949152542373
$ws.Range("S9").Value = 0.4293785310734463
$ws.Range("B10").Value = 0.1134099616858237
$ws.Range("D10").Value = 0.02681992337164751
$ws.Range("E10").Value = 0.001532567049808429
$ws.Range("F10").Value = 0.05747126436781609
$ws.Range("J10").Value = 0.1210727969348659
$ws.Range("O10").Value = 0.009195402298850575
$ws.Range("Q10").Value = 0.2245210727969349
$ws.Range("R10").Value = 0.07892720306513409
$ws.Range("S10").Value = 0.3670498084291188
$ws.Range("G11").Value = 0.125
$ws.Range("J11").Value = 0.06510416666666667
$ws.Range("K11").Value = 0.1640625
$ws.Range("L11").Value = 0.6197916666666666
$ws.Range("S11").Value = 0.02604166666666667
$ws.Range("G12").Value = 0.7708333333333334
$ws.Range("J12").Value = 0.1958333333333333
$ws.Range("K12").Value = 0.008333333333333333
$ws.Range("L12").Value = 0.008333333333333333
$ws.Range("S12").Value = 0.01666666666666667
$ws.Range("G13").Value = 0.5434782608695652
$ws.Range("J13").Value = 0.391304347826087
$ws.Range("S13").Value = 0.06521739130434782
$ws.Range("F15").Value = 0.0202020202020202
$ws.Range("H15").Value = 0.1666666666666667
$ws.Range("I15").Value = 0.06565656565656566
$ws.Range("J15").Value = 0.3181818181818182
$ws.Range("K15").Value = 0.09595959595959595
$ws.Range("O15").Value = 0.03535353535353535
$ws.Range("S15").Value = 0.297979797979798
$ws.Range("F16").Value = 0.009852216748768473
$ws.Range("H16").Value = 0.2266009852216749
$ws.Range("I16").Value = 0.06403940886699508
$ws.Range("J16").Value = 0.3793103448275862
$ws.Range("K16").Value = 0.1133004926108374
$ws.Range("M16").Value = 0.01477832512315271
$ws.Range("O16").Value = 0.0541871921182266
$ws.Range("S16").Value = 0.1379310344827586
$ws.Range("F17").Value = 0.02569169960474308
$ws.Range("H17").Value = 0.183794466403162
$ws.Range("I17").Value = 0.0691699604743083
$ws.Range("J17").Value = 0.3952569169960474
$ws.Range("K17").Value = 0.1225296442687747
$ws.Range("M17").Value = 0.009881422924901186
$ws.Range("N17").Value = 0.001976284584980237
$ws.Range("O17").Value = 0.05533596837944664
$ws.Range("S17").Value = 0.1363636363636364
$ws.Range("F18").Value = 0.0187793427230047
$ws.Range("H18").Value = 0.2065727699530517
$ws.Range("I18").Value = 0.07981220657276995
$ws.Range("J18").Value = 0.4037558685446009
$ws.Range("K18").Value = 0.1126760563380282
$ws.Range("M18").Value = 0.02347417840375587
$ws.Range("O18").Value = 0.02816901408450704
$ws.Range("S18").Value = 0.1267605633802817
$ws.Range("F19").Value = 0.01268498942917548
$ws.Range("H19").Value = 0.2367864693446089
$ws.Range("I19").Value = 0.07047216349541931
$ws.Range("J19").Value = 0.3192389006342495
$ws.Range("K19").Value = 0.1303735024665257
$ws.Range("M19").Value = 0.02325581395348837
$ws.Range("O19").Value = 0.06342494714587738
$ws.Range("S19").Value = 0.1437632135306554
